# Reassign the per-observation data among rows 31-57 of the Artfynd sheet.
# The row's "identity" columns (A id, Q/R coords, Z/AB times together with the
# species/record columns B,D,E,F,G,H,M) are permuted across rows while the
# shared columns (location/municipality/observer, etc.) stay put since they
# are identical for every row in this block anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values travel together with a given observation record.
$cols = @("A","B","D","E","F","G","H","M","Q","R","Z","AB")

# Destination row -> source row (source = where that data currently lives).
$rowMap = @{
    31 = 33
    33 = 31
    34 = 35
    35 = 36
    36 = 37
    37 = 34
    39 = 41
    41 = 40
    40 = 42
    42 = 39
    43 = 44
    44 = 43
    45 = 46
    46 = 45
    48 = 49
    49 = 48
    52 = 54
    54 = 53
    53 = 52
    55 = 57
    57 = 55
}

# First, snapshot the current ("before") values of every cell we might need,
# since several destinations read from rows that are themselves overwritten.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the permuted values back out.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $val = $srcVals[$c]
        if ($val -eq $null) {
            $ws.Range("$c$destRow").Value = ""
        } else {
            $ws.Range("$c$destRow").Value = $val
        }
    }
}
